$wb = $excel.ActiveWorkbook

# Snapshot the sheets that exist before we add the new "gather" sheet.
$sourceSheets = @()
foreach ($s in $wb.Worksheets) {
    $sourceSheets += $s
}

# Add the new sheet after the last existing sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "gather_2"

$colLetters = @("A", "B", "C", "D")

for ($i = 0; $i -lt $sourceSheets.Count; $i++) {
    $src = $sourceSheets[$i]
    $colLetter = $colLetters[$i]
    $destCol = $i + 1

    # Header = source sheet name, written as text (not a number) even
    # though the sheet names look numeric.
    $headerCell = $newSheet.Range($colLetter + "1")
    $headerCell.NumberFormat = "@"
    $headerCell.Value2 = $src.Name

    # Column B on each source sheet holds the "Time" values we gather.
    $lastRow = $src.Cells.Item($src.Rows.Count, 2).End(-4162).Row

    $destRow = 2
    for ($r = 2; $r -le $lastRow; $r++) {
        $val = $src.Cells.Item($r, 2).Value2
        $newSheet.Cells.Item($destRow, $destCol).Value2 = $val
        $destRow++
    }
}
